$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '36.416.55'
$ws.Range('E2').Value = '  +2.81%  '
# Row 3
$ws.Range('D3').Value = '2.011.35'
$ws.Range('E3').Value = '  +6.09%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').Value = '''245.23'
$ws.Range('E5').Value = '  -0.40%  '
# Row 6
$ws.Range('D6').Value = '''0.660'
$ws.Range('E6').Value = '  -4.71%  '
# Row 7
$ws.Range('E7').Value = '  +0.00%  '
# Row 8
$ws.Range('D8').Value = '''45.01'
$ws.Range('E8').Value = '  +4.52%  '
# Row 9
$ws.Range('D9').Value = '''60.70'
$ws.Range('E9').Value = '  +7.98%  '
# Row 10
$ws.Range('D10').Value = '''0.370'
$ws.Range('E10').Value = '  +3.56%  '
# Row 11
$ws.Range('D11').Value = '''0.0714'
$ws.Range('E11').Value = '  -5.63%  '
# Row 12
$ws.Range('E12').Value = '  +0.18%  '
# Row 13
$ws.Range('D13').Value = '''14.54'
$ws.Range('E13').Value = '  +2.60%  '
# Row 14
$ws.Range('D14').Value = '2.308.57'
$ws.Range('E14').Value = '  +6.27%  '
# Row 15
$ws.Range('D15').Value = '''0.809'
$ws.Range('E15').Value = '  +1.48%  '
# Row 16
$ws.Range('D16').Value = '2.021.80'
$ws.Range('E16').Value = '  +7.95%  '
# Row 17
$ws.Range('E17').Value = '  -2.22%  '
# Row 18
$ws.Range('D18').Value = '36.270.15'
$ws.Range('E18').Value = '  +2.39%  '
# Row 19
$ws.Range('D19').Value = '''71.33'
$ws.Range('E19').Value = '  -3.14%  '
# Row 20
$ws.Range('E20').Value = '  -1.87%  '
# Row 21
$ws.Range('D21').Value = '''12.85'
$ws.Range('E21').Value = '  -1.13%  '
# Row 22
$ws.Range('D22').Value = '''236.06'
$ws.Range('E22').Value = '  -3.43%  '
# Row 23
$ws.Range('E23').Value = '  -6.39%  '
# Row 24
$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  -0.12%  '
# Row 25
$ws.Range('D25').Value = '''2.43'
$ws.Range('E25').Value = '  -9.50%  '
# Row 26
$ws.Range('D26').Value = '''163.94'
$ws.Range('E26').Value = '  -1.69%  '
# Row 27
$ws.Range('D27').Value = '''19.61'
$ws.Range('E27').Value = '  +7.10%  '
# Row 28
$ws.Range('D28').Value = '''8.56'
$ws.Range('E28').Value = '  -0.55%  '
# Row 29
$ws.Range('E29').Value = '  -10.70%  '
# Row 30
$ws.Range('E30').Value = '  -4.77%  '
# Row 31
$ws.Range('D31').Value = '''22.51'
$ws.Range('E31').Value = '  +63.38%  '
# Row 32
$ws.Range('D32').Value = '''4.40'
$ws.Range('E32').Value = '  +0.90%  '
# Row 33
$ws.Range('D33').Value = '''0.0586'
$ws.Range('E33').Value = '  -2.83%  '
# Row 34
$ws.Range('E34').Value = '  +0.01%  '
# Row 35
$ws.Range('D35').Value = '''1.87'
$ws.Range('E35').Value = '  -0.39%  '
# Row 36
$ws.Range('D36').Value = '''3.99'
$ws.Range('E36').Value = '  -5.96%  '
# Row 37
$ws.Range('D37').Value = '''2.17'
$ws.Range('E37').Value = '  +11.08%  '
# Row 38
$ws.Range('D38').Value = '''0.0811'
$ws.Range('E38').Value = '  +9.98%  '
# Row 39
$ws.Range('D39').Value = '''0.850'
$ws.Range('E39').Value = '  -0.57%  '
# Row 40
$ws.Range('D40').Value = '''1.34'
$ws.Range('E40').Value = '  -9.75%  '
# Row 41
$ws.Range('D41').Value = '''0.0217'
$ws.Range('E41').Value = '  -3.62%  '
# Row 42
$ws.Range('D42').Value = '''95.91'
$ws.Range('E42').Value = '  -3.14%  '
# Row 43
$ws.Range('D43').Value = '''1.11'
$ws.Range('E43').Value = '  +2.44%  '
# Row 44
$ws.Range('E44').Value = '  +14.94%  '
# Row 45
$ws.Range('D45').Value = '''16.00'
$ws.Range('E45').Value = '  -5.81%  '
# Row 46
$ws.Range('D46').Value = '1.314.38'
$ws.Range('E46').Value = '  -0.85%  '
# Row 47
$ws.Range('E47').Value = '  +0.83%  '
# Row 48
$ws.Range('D48').Value = '''2.76'
$ws.Range('E48').Value = '  +1.01%  '
# Row 49
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.199.86'
$ws.Range('E49').Value = '  +6.11%  '
# Row 50
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '''2.21'
$ws.Range('E50').Value = '  -6.87%  '
# Row 51
$ws.Range('E51').Value = '  +15.56%  '
